# Update column F (dSF) values for specific rows, per repull/push of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 1
    6  = -3
    12 = -1
    13 = 0
    17 = -1
    18 = 1
    23 = 6
    35 = -2
    36 = -3
    38 = 5
    39 = 1
    42 = 3
    43 = 6
    44 = -4
    50 = -4
    51 = -3
    53 = -5
    55 = -3
    56 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
